$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CHEMICAL PROD)
$ws.Range("B2").Value = [double]"-46.63591686726431"
$ws.Range("H2").Value = [double]"-0.08093757796183354"

# Row 3 (COFFEE ESTATE)
$ws.Range("B3").Value = [double]"2.795794920751327e-05"
$ws.Range("H3").Value = [double]"0.002901338853112279"

# Row 4 (COOPERATIVES)
$ws.Range("B4").Value = [double]"4.941826115461026"
$ws.Range("H4").Value = [double]"0.1508299285860346"

# Row 5 (ELECTRICITY PROD)
$ws.Range("B5").Value = [double]"-0.6628334176880344"
$ws.Range("D5").Value = [double]"-9.378236803240725"
$ws.Range("E5").Value = [double]"-2.162528331144131"
$ws.Range("H5").Value = [double]"-4.723281139569735"
$ws.Range("J5").Value = [double]"-0.007281241271158478"

# Row 6 (FERTILIZERS PROD)
$ws.Range("B6").Value = [double]"-45.23583239729487"
$ws.Range("H6").Value = [double]"-0.07850770301670451"

# Row 7 (INFORMAL)
$ws.Range("B7").Value = [double]"-122.3768405538212"
$ws.Range("H7").Value = [double]"-3.735074786134192"

# Row 8 (MANUFACTURING)
$ws.Range("B8").Value = [double]"-0.006652535451003416"
$ws.Range("C8").Value = [double]"-1.627611016474475"
$ws.Range("H8").Value = [double]"-4.122085544087895"

# Row 9 (PETROLEUM PROD)
$ws.Range("B9").Value = [double]"0.007174824431658067"
$ws.Range("H9").Value = [double]"-0.7123601044104362"

# Row 10 (PRIMARY)
$ws.Range("B10").Value = [double]"-0.3303854311901659"
$ws.Range("D10").Value = [double]"-4.652951300662608"
$ws.Range("E10").Value = [double]"-1.072924482992676"
$ws.Range("H10").Value = [double]"-2.436277699363927"
$ws.Range("J10").Value = [double]"-0.003612539426587347"

# Row 11 (SERVICES)
$ws.Range("B11").Value = [double]"-0.005273460015359888"
$ws.Range("H11").Value = [double]"-3.775628661343944"
